$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.686.88'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '3.803.04'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '596.81'
$ws.Range('D6').Value = '167.24'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.57%  '
$ws.Range('E9').Value = '  +1.23%  '
$ws.Range('D11').Value = '0.450'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').Value = '35.90'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '4.442.28'
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').Value = '3.804.27'
$ws.Range('E15').Value = '  +0.87%  '
$ws.Range('D16').Value = '18.61'
$ws.Range('E16').Value = '  +5.03%  '
$ws.Range('D17').Value = '67.709.59'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('E18').Value = '  +2.17%  '
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').Value = '461.35'
$ws.Range('E20').Value = '  +1.31%  '
$ws.Range('D21').Value = '9.94'
$ws.Range('E21').Value = '  -2.25%  '
$ws.Range('E22').Value = '  +0.90%  '
$ws.Range('E23').Value = '  +0.50%  '
$ws.Range('D24').Value = '83.42'
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('E25').Value = '  +2.40%  '
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').Value = '3.944.44'
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').Value = '2.23'
$ws.Range('E31').Value = '  +2.72%  '
$ws.Range('D32').Value = '7.30'
$ws.Range('E32').Value = '  +1.48%  '
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E35').Value = '  -0.98%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.746.55'
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '0.0999'
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').Value = '3.38'
$ws.Range('E38').Value = '  +2.59%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '0.138'
$ws.Range('E39').Value = '  +0.35%  '
$ws.Range('B40').Value = 'Mantle'
$ws.Range('C40').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  +0.79%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '5.77'
$ws.Range('E41').Value = '  +0.83%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '48.16'
$ws.Range('E44').Value = '  +3.00%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '0.300'
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('B46').Value = 'Arweave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D46').Value = '43.04'
$ws.Range('E46').Value = '  -1.33%  '
$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D47').Value = '8.33'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = '27.18'
$ws.Range('E48').Value = '  +8.24%  '
$ws.Range('D49').Value = '147.36'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').Value = '1.36'
$ws.Range('E50').Value = '  +11.63%  '
$ws.Range('B51').Value = 'Bittensor'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D51').Value = '394.07'
$ws.Range('E51').Value = '  +1.36%  '
